$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix punctuation in proveedor/contratista name fields (stray comma -> period) ---
$ws.Range("E59").Value  = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E70").Value  = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
$ws.Range("F70").Value  = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
$ws.Range("E77").Value  = "RICCOTTI. MARIANA EDITH"
$ws.Range("E108").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"

# --- Fix "Importe" (column H) values scraped with Spanish-style grouping, e.g. "1.234,56" ---
# --- Re-enter as plain-text numbers with a dot decimal separator and no thousands separator ---
# --- (format the column as Text first so Excel keeps these as strings, not numbers) ---
$ws.Range("H2:H148").NumberFormat = "@"

$importeRows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76, 77, 78, 79, 80, 81, 82, 83, 84, 85, 86, 87, 88, 89, 90, 91, 92, 93, 94, 95, 96, 97, 98, 99, 100, 101, 102, 103, 104, 105, 106, 107, 108, 109, 110, 111, 112, 113, 114, 115, 116, 117, 118, 119, 120, 121, 122, 123, 124, 125, 126, 127, 128, 129, 130, 131, 132, 133, 134, 135, 136, 137, 138, 139, 140, 141, 142, 143, 144, 145, 146, 147, 148)
$importeVals = @("2918.50", "22870.00", "91935.52", "298286.79", "2378.00", "2550.00", "2400.86", "479.76", "8961.58", "687.50", "14.40", "7899.61", "2320.00", "532.15", "60.00", "237.50", "2960.00", "217.06", "22992.54", "247.50", "1643.56", "27.70", "647.00", "1725.00", "3001.00", "10936.65", "3190.23", "645.52", "4356.27", "2390.00", "7411.97", "15977.12", "750.00", "80.00", "136.00", "13342.34", "2283.00", "7260.00", "10544.79", "575.50", "25.12", "420.00", "4500.00", "1090.00", "18.05", "7595.00", "11928.00", "2181.00", "2610.00", "3610.64", "750.00", "152.00", "3760.00", "959.00", "1155.00", "250.00", "20256.00", "858.00", "18254.00", "88.00", "154623.41", "378.48", "96.00", "420.00", "4404.50", "2400.00", "1140.90", "313.00", "7087.00", "6117.95", "8.75", "183.81", "1339.00", "60.00", "1953.20", "10000.00", "1504.00", "720.00", "798.00", "1772.00", "1793.22", "600.00", "1885.00", "871.69", "1710.00", "90720.00", "2670.00", "2000.00", "250.00", "1200.00", "2000.00", "700.00", "250.00", "600.00", "4000.00", "21326.51", "500.00", "950.00", "650.00", "100.00", "2450.00", "1500.00", "200.00", "240.00", "19160.00", "556.15", "20.00", "85.00", "1350.00", "7200.00", "1334.03", "1575.00", "9705.00", "64.70", "214.64", "77.00", "57.80", "2634.00", "370.00", "543.71", "201.20", "273.29", "2996.00", "518.00", "878.00", "550.20", "953.57", "2668.14", "877.86", "13200.00", "490.00", "692.92", "790560.12", "1100.00", "1260.00", "402800.00", "34000.00", "407904.50", "20000.00", "427954.00", "417968.00", "414000.00", "409048.00", "380000.00", "411164.69", "605.00", "140.00")

for ($i = 0; $i -lt $importeRows.Count; $i++) {
    $ws.Cells.Item($importeRows[$i], 8).Value = $importeVals[$i]
}
